# Horarios Línea 141 - update scrape results (08:11:27 -> 08:29:19 run)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": header text, row swap, new rows inserted mid-table, and
# several brand-new rows appended from the latest scrape.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 08:29:19"
$ws1.Cells.Item(3,1).Value = "Total filas: 96"

# Old rows 52/53 come out swapped in the refreshed scrape.
$ws1.Cells.Item(52,1).Value = "07:36:59"
$ws1.Cells.Item(52,2).Value = "08:02"
$ws1.Cells.Item(52,3).Value = "17_ROMERO"
$ws1.Cells.Item(52,4).Value = 26
$ws1.Cells.Item(52,5).Value = "LP1912"

$ws1.Cells.Item(53,1).Value = "06:52:52"
$ws1.Cells.Item(53,2).Value = "08:02"
$ws1.Cells.Item(53,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(53,4).Value = 70
$ws1.Cells.Item(53,5).Value = "LP1912"

# New row inserted at row 64, pushing the former rows 64-93 down by one.
$ws1.Rows.Item(64).Insert()
$ws1.Cells.Item(64,1).Value = "08:29:19"
$ws1.Cells.Item(64,2).Value = "08:29"
$ws1.Cells.Item(64,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(64,4).Value = 0
$ws1.Cells.Item(64,5).Value = "LP1912"

# New row inserted at row 76 (post-shift numbering).
$ws1.Rows.Item(76).Insert()
$ws1.Cells.Item(76,1).Value = "08:29:19"
$ws1.Cells.Item(76,2).Value = "08:54"
$ws1.Cells.Item(76,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(76,4).Value = 25
$ws1.Cells.Item(76,5).Value = "LP1912"

# New row inserted at row 81 (post-shift numbering).
$ws1.Rows.Item(81).Insert()
$ws1.Cells.Item(81,1).Value = "08:29:19"
$ws1.Cells.Item(81,2).Value = "09:03"
$ws1.Cells.Item(81,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(81,4).Value = 34
$ws1.Cells.Item(81,5).Value = "LP1912"

# Five brand-new rows appended at the bottom (97-101).
$ws1.Cells.Item(97,1).Value = "08:29:19"
$ws1.Cells.Item(97,2).Value = "10:11"
$ws1.Cells.Item(97,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(97,4).Value = 102
$ws1.Cells.Item(97,5).Value = "LP1912"

$ws1.Cells.Item(98,1).Value = "08:29:19"
$ws1.Cells.Item(98,2).Value = "10:12"
$ws1.Cells.Item(98,3).Value = "10_OLMOS"
$ws1.Cells.Item(98,4).Value = 103
$ws1.Cells.Item(98,5).Value = "LP1912"

$ws1.Cells.Item(99,1).Value = "08:29:19"
$ws1.Cells.Item(99,2).Value = "10:14"
$ws1.Cells.Item(99,3).Value = "10_OLMOS"
$ws1.Cells.Item(99,4).Value = 105
$ws1.Cells.Item(99,5).Value = "LP1912"

$ws1.Cells.Item(100,1).Value = "08:29:19"
$ws1.Cells.Item(100,2).Value = "10:15"
$ws1.Cells.Item(100,3).Value = "17_ROMERO"
$ws1.Cells.Item(100,4).Value = 106
$ws1.Cells.Item(100,5).Value = "LP1912"

$ws1.Cells.Item(101,1).Value = "08:29:19"
$ws1.Cells.Item(101,2).Value = "10:26"
$ws1.Cells.Item(101,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(101,4).Value = 117
$ws1.Cells.Item(101,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": only the "last updated" timestamp changes.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 08:29:19"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": header text plus one brand-new row appended at the end.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 08:29:19"
$ws3.Cells.Item(3,1).Value = "Total filas: 12"

$ws3.Cells.Item(17,1).Value = "08:29:19"
$ws3.Cells.Item(17,2).Value = "10:23"
$ws3.Cells.Item(17,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(17,4).Value = 114
$ws3.Cells.Item(17,5).Value = "L6173"
